$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new year column (2021) and its data value (515),
# matching the formatting of the adjacent column N
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 515

# Update selection / view state to match the saved workbook
$ws.Range("P12").Select()

$wb.Save()
